# Regenerate merged AHB files
# - Rename the "_old" / "_new" suffixed header columns to "_FV2410" / "_FV2504"
# - Add a structured table (ListObject) over the full data range
# - Freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 holds the column headers; columns A-J end in "_old" and columns L-U end
# in "_new" (column K is the literal "diff" column and is left untouched).
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $header = [string]$cell.Value2
    if ($header -like "*_old") {
        $cell.Value = $header.Substring(0, $header.Length - 4) + "_FV2410"
    } elseif ($header -like "*_new") {
        $cell.Value = $header.Substring(0, $header.Length - 4) + "_FV2504"
    }
}

# Turn the data range into a structured table
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U78"), $null, 1)
$lo.Name = "Table1"

# Freeze the header row
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
